# Generate Report for Handback
#
# - Overview sheet: zh-cn / de-de status cells flip from "Ready for handoff"
#   to "Handed back: in sync with en-US" (widens the status columns to fit).
# - zh-cn / de-de sheets: fill in the "Latest Target File" (col I) and
#   "Latest Handback File" (col J) for both data rows, add hyperlinks on the
#   new target-file cells (mirroring the existing source-file hyperlinks),
#   stamp the "Latest Handback DateTime" (col K), and widen the Status /
#   Target-File / Handback-File columns to fit the new content.

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dff813a533eb67f336a7b270fdc4f50ca97682/e2e/"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status cells for both rows.
# ---------------------------------------------------------------------
$overview.Cells.Item(2, 5).Value = $newStatus   # E2 zh-cn status
$overview.Cells.Item(2, 6).Value = $newStatus   # F2 de-de status
$overview.Cells.Item(3, 5).Value = $newStatus   # E3 zh-cn status
$overview.Cells.Item(3, 6).Value = $newStatus   # F3 de-de status

# Widen columns E & F on the Overview sheet to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet (row 2 -> 0186399d..., row 3 -> 75442726...)
# ---------------------------------------------------------------------
$zhcn.Cells.Item(2, 9).Value  = "0186399d-6535-49f4-a04b-eaba33fdb113.md"
$zhcn.Cells.Item(2, 10).Value = "0186399d-6535-49f4-a04b-eaba33fdb113.283cec73cb65d45d1ab87d0764375ffb043fe3ca.zh-cn.xlf"
$zhcn.Cells.Item(2, 11).Value = "2016-08-21 11:05:07"

$zhcn.Cells.Item(3, 9).Value  = "75442726-3f64-4f81-94b5-f3fdcef57e23.md"
$zhcn.Cells.Item(3, 10).Value = "75442726-3f64-4f81-94b5-f3fdcef57e23.a49100b713f4387c607885921aa3fec1ae93ed1e.zh-cn.xlf"
$zhcn.Cells.Item(3, 11).Value = "2016-08-21 11:05:07"

# Style the new "Latest Target File" cells like the existing hyperlink cells.
$zhcn.Cells.Item(2, 9).Style = "HyperLink"
$zhcn.Cells.Item(3, 9).Style = "HyperLink"

# Column widths: Status (C) and the two new file columns (I, J).
$zhcn.Columns.Item(3).ColumnWidth  = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth  = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet (row 2 -> 0186399d..., row 3 -> 75442726...)
# ---------------------------------------------------------------------
$dede.Cells.Item(2, 9).Value  = "0186399d-6535-49f4-a04b-eaba33fdb113.md"
$dede.Cells.Item(2, 10).Value = "0186399d-6535-49f4-a04b-eaba33fdb113.283cec73cb65d45d1ab87d0764375ffb043fe3ca.de-de.xlf"
$dede.Cells.Item(2, 11).Value = "2016-08-21 11:05:14"

$dede.Cells.Item(3, 9).Value  = "75442726-3f64-4f81-94b5-f3fdcef57e23.md"
$dede.Cells.Item(3, 10).Value = "75442726-3f64-4f81-94b5-f3fdcef57e23.a49100b713f4387c607885921aa3fec1ae93ed1e.de-de.xlf"
$dede.Cells.Item(3, 11).Value = "2016-08-21 11:05:14"

$dede.Cells.Item(2, 9).Style = "HyperLink"
$dede.Cells.Item(3, 9).Style = "HyperLink"

$dede.Columns.Item(3).ColumnWidth  = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth  = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# Hyperlinks: rebuild each sheet's hyperlink collection so the source-file
# link (col A) and the new target-file link (col I) both exist, in row
# order, for row 2 then row 3.
# ---------------------------------------------------------------------
foreach ($ws in @($zhcn, $dede)) {
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Cells.Item(2, 1), ($ghBase + "0186399d-6535-49f4-a04b-eaba33fdb113.md"), [Type]::Missing, [Type]::Missing, "0186399d-6535-49f4-a04b-eaba33fdb113.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(2, 9), ($ghBase + "0186399d-6535-49f4-a04b-eaba33fdb113.md"), [Type]::Missing, [Type]::Missing, "0186399d-6535-49f4-a04b-eaba33fdb113.md")

    $ws.Hyperlinks.Add($ws.Cells.Item(3, 1), ($ghBase + "75442726-3f64-4f81-94b5-f3fdcef57e23.md"), [Type]::Missing, [Type]::Missing, "75442726-3f64-4f81-94b5-f3fdcef57e23.md")
    $ws.Hyperlinks.Add($ws.Cells.Item(3, 9), ($ghBase + "75442726-3f64-4f81-94b5-f3fdcef57e23.md"), [Type]::Missing, [Type]::Missing, "75442726-3f64-4f81-94b5-f3fdcef57e23.md")
}
